$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: corrected light/dark timing values ---
$ws.Range("E10").Value = 45
$ws.Range("F10").Value = 16
$ws.Range("I10").Value = 6
$ws.Range("AD10").Value = 23

# --- Row 11: corrected light/dark timing values ---
$ws.Range("E11").Value = 45
$ws.Range("F11").Value = 16
$ws.Range("I11").Value = 56

# --- Row 12: corrected light/dark timing values ---
$ws.Range("E12").Value = 45
$ws.Range("F12").Value = 16
$ws.Range("I12").Value = 56

# Row 12: fill in sample date/time and 1st-3rd tube change info (flask K)
$ws.Range("P12").Value = 44257
$ws.Range("Q12").Value = 2021
$ws.Range("R12").Value = 3
$ws.Range("S12").Value = 2
$ws.Range("T12").Value = 8
$ws.Range("U12").Value = 38
$ws.Range("V12").Value = 2021
$ws.Range("W12").Value = 3
$ws.Range("X12").Value = 2
$ws.Range("Y12").Value = 15
$ws.Range("Z12").Value = 45
$ws.Range("AA12").Value = 2021
$ws.Range("AB12").Value = 3
$ws.Range("AC12").Value = 2
$ws.Range("AD12").Value = 23
$ws.Range("AE12").Value = 18
$ws.Range("AF12").Value = 2021
$ws.Range("AG12").Value = 3
$ws.Range("AH12").Value = 3
$ws.Range("AI12").Value = 8
$ws.Range("AJ12").Value = 31
$ws.Range("AA12:AJ12").Font.Color = 0

# --- Row 13: corrected/added light-dark duration values (flask L) ---
$ws.Range("E13").Value = 45
$ws.Range("F13").Value = 21
$ws.Range("G13").Value = 4
$ws.Range("H13").Value = 12
$ws.Range("I13").Value = 19

# Row 13: fill in sample date/time and 1st-4th tube change info
$ws.Range("P13").Value = 44258
$ws.Range("Q13").Value = 2021
$ws.Range("R13").Value = 3
$ws.Range("S13").Value = 3
$ws.Range("T13").Value = 8
$ws.Range("U13").Value = 38
$ws.Range("V13").Value = 2021
$ws.Range("W13").Value = 3
$ws.Range("X13").Value = 3
$ws.Range("Y13").Value = 17
$ws.Range("Z13").Value = 31
$ws.Range("AA13").Value = 2021
$ws.Range("AB13").Value = 3
$ws.Range("AC13").Value = 3
$ws.Range("AD13").Value = 19
$ws.Range("AE13").Value = 38
$ws.Range("AF13").Value = 2021
$ws.Range("AG13").Value = 3
$ws.Range("AH13").Value = 3
$ws.Range("AI13").Value = 23
$ws.Range("AJ13").Value = 35
$ws.Range("AK13").Value = 2021
$ws.Range("AL13").Value = 3
$ws.Range("AM13").Value = 4
$ws.Range("AN13").Value = 9
$ws.Range("AO13").Value = 38
$ws.Range("Q13:U13").Font.Color = 0

# --- Update the active selection to reflect latest data entry position ---
$ws.Range("I19").Select() | Out-Null
